$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing last row (row 75) values
$ws.Range("B75").Value = 20.9
$ws.Range("C75").Value = 0.8
$ws.Range("D75").Value = -9.5

# Add new row 76 with the next quarterly period.
# Force text formatting first so the quarter label (which looks like a
# date) is stored as literal text instead of being auto-converted to a
# date serial value, then restore the default (unformatted) style so the
# new cell matches the rest of column A.
$ws.Range("A76").NumberFormatLocal = "@"
$ws.Range("A76").Value = "01-04-2021"
$ws.Range("A76").Style = $ws.Range("A75").Style
$ws.Range("B76").Value = 19.5
$ws.Range("C76").Value = -1.1
$ws.Range("D76").Value = -5.6
